# Generate Report for Handback
# Update the "Latest HO Xliff Generate Date" / handoff / handback datetime
# stamps that are refreshed whenever the handback report is (re)generated.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for b01085d3...md row
$wsOverview.Range("G3").Value = "2016-08-21 12:51:19"

# zh-cn sheet, row for b01085d3...: Correspond Handoff Datetime / Correspond Handback DateTime
$wsZhCn.Range("H3").Value = "2016-08-21 12:51:15"
$wsZhCn.Range("K3").Value = "2016-08-21 12:51:29"

# de-de sheet, row for b01085d3...: Correspond Handoff Datetime (shares the same
# original value/shared-string slot as Overview!G3) and Correspond Handback DateTime
$wsDeDe.Range("H3").Value = "2016-08-21 12:51:19"
$wsDeDe.Range("K3").Value = "2016-08-21 12:51:36"
